# This script updates betting-odds values on Sheet1 to reflect the
# latest FlashScore data refresh, as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 4
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 2.08
$ws.Range("R4").Value = 1.73

# Row 6
$ws.Range("O6").Value = 1.17
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.57
$ws.Range("R6").Value = 2.35

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 4.5
$ws.Range("L7").Value = 4.75
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 21
$ws.Range("AM7").Value = 301
$ws.Range("AP7").Value = 23
$ws.Range("AX7").Value = 23
$ws.Range("BA7").Value = 101

# Row 9
$ws.Range("Q9").Value = 1.95
$ws.Range("R9").Value = 1.95

# Row 14
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("Q14").Value = 1.92
$ws.Range("R14").Value = 1.98

# Row 16
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 7

# Row 17
$ws.Range("G17").Value = 2
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 2.75
$ws.Range("L17").Value = 4.75
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 8
$ws.Range("AD17").Value = 6.5
$ws.Range("AN17").Value = 3.75
